$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '89.874.73'
$ws.Cells.Item(2, 5).Value = '  -1.46%  '
$ws.Cells.Item(3, 4).Value = '3.102.96'
$ws.Cells.Item(3, 5).Value = '  -2.52%  '
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '0.995'
$ws.Cells.Item(4, 5).Value = '  -0.47%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '214.30'
$ws.Cells.Item(5, 5).Value = '  -0.96%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '619.66'
$ws.Cells.Item(6, 5).Value = '  -2.33%  '
$ws.Cells.Item(7, 5).Value = '  -5.32%  '
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.868'
$ws.Cells.Item(8, 5).Value = '  +21.18%  '
$ws.Cells.Item(9, 5).Value = '  +0.04%  '
$ws.Cells.Item(10, 4).Value = '3.099.60'
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.634'
$ws.Cells.Item(11, 5).Value = '  +11.84%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.189'
$ws.Cells.Item(12, 5).Value = '  +4.25%  '
$ws.Cells.Item(13, 5).Value = '  -5.37%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '5.32'
$ws.Cells.Item(14, 5).Value = '  -0.02%  '
$ws.Cells.Item(15, 4).Value = '89.665.18'
$ws.Cells.Item(15, 5).Value = '  -1.25%  '
$ws.Cells.Item(16, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(16, 4).Value = '3.691.30'
$ws.Cells.Item(16, 5).Value = '  -2.12%  '
$ws.Cells.Item(17, 2).Value = 'Avalanche'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '32.37'
$ws.Cells.Item(17, 5).Value = '  -0.40%  '
$ws.Cells.Item(18, 4).Value = '3.113.40'
$ws.Cells.Item(18, 5).Value = '  -1.96%  '
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '3.42'
$ws.Cells.Item(19, 5).Value = '  +2.91%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '0.0000219'
$ws.Cells.Item(20, 5).Value = '  +2.28%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '13.45'
$ws.Cells.Item(21, 5).Value = '  +1.30%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '427.48'
$ws.Cells.Item(22, 5).Value = '  -1.60%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '8.26'
$ws.Cells.Item(23, 5).Value = '  -2.47%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '4.97'
$ws.Cells.Item(24, 5).Value = '  -0.26%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '5.44'
$ws.Cells.Item(25, 5).Value = '  +3.89%  '
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '12.04'
$ws.Cells.Item(26, 5).Value = '  +3.54%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '83.92'
$ws.Cells.Item(27, 5).Value = '  +4.04%  '
$ws.Cells.Item(28, 4).Value = '3.299.76'
$ws.Cells.Item(28, 5).Value = '  -1.59%  '
$ws.Cells.Item(29, 5).Value = '  +0.04%  '
$ws.Cells.Item(30, 5).Value = '  +8.70%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '0.162'
$ws.Cells.Item(31, 5).Value = '  +0.08%  '
$ws.Cells.Item(32, 5).Value = '  -2.20%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '511.82'
$ws.Cells.Item(33, 5).Value = '  -1.00%  '
$ws.Cells.Item(34, 5).Value = '  -7.55%  '
$ws.Cells.Item(35, 5).Value = '  -3.45%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '1.25'
$ws.Cells.Item(36, 5).Value = '  -4.01%  '
$ws.Cells.Item(37, 5).Value = '  -4.63%  '
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '22.56'
$ws.Cells.Item(38, 5).Value = '  +0.98%  '
$ws.Cells.Item(39, 5).Value = '  -0.46%  '
$ws.Cells.Item(40, 5).Value = '  +3.61%  '
$ws.Cells.Item(41, 5).Value = '  +0.18%  '
$ws.Cells.Item(42, 5).Value = '  -0.02%  '
$ws.Cells.Item(43, 2).Value = 'Stellar'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '0.142'
$ws.Cells.Item(43, 5).Value = '  +13.61%  '
$ws.Cells.Item(44, 2).Value = 'PolygonEcosystemToken'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '0.368'
$ws.Cells.Item(44, 5).Value = '  -0.30%  '
$ws.Cells.Item(45, 2).Value = 'Stacks'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '1.85'
$ws.Cells.Item(45, 5).Value = '  -3.43%  '
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '145.94'
$ws.Cells.Item(46, 5).Value = '  -0.85%  '
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '0.0703'
$ws.Cells.Item(47, 5).Value = '  +14.36%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '43.63'
$ws.Cells.Item(48, 5).Value = '  -1.21%  '
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '1.22'
$ws.Cells.Item(49, 5).Value = '  +1.37%  '
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '159.47'
$ws.Cells.Item(50, 5).Value = '  -5.98%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '0.709'
$ws.Cells.Item(51, 5).Value = '  -4.15%  '
